$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 162, shifting existing rows 162-239 down to 163-240
$ws.Rows.Item(162).Insert()

# Populate the newly inserted row 162 with the new data record
$ws.Range("A162").Value = 9
$ws.Range("B162").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C162").Value = "Metropolitana"
$ws.Range("D162").Value = 44510
$ws.Range("E162").Value = 13
$ws.Range("F162").Value = 100112044
$ws.Range("G162").Value = "Perejil"
$ws.Range("H162").Value = "Sin especificar"
$ws.Range("I162").Value = "Primera"
$ws.Range("J162").Value = 97
$ws.Range("K162").Value = 12000
$ws.Range("L162").Value = 13000
$ws.Range("M162").Value = 12495
$ws.Range("N162").Value = "$/docena de atados"
$ws.Range("O162").Value = "Región Metropolitana"
$ws.Range("P162").Value = 4165
$ws.Range("Q162").Value = 3
$ws.Range("R162").Value = "Hortaliza"
